$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume 1h (E) columns with the latest scraped values.
# Cells whose new price text is an unambiguous decimal number (e.g. "0.9986")
# are explicitly formatted as Text first so Excel keeps them as strings instead
# of silently converting them to numeric cells (values with multiple "." like
# "29.147.99" are already unambiguous text and do not need this).

$ws.Range('D2').Value = '29.147.99'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '1.835.14'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9986'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E5').Value = '  -1.96%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6652'
$ws.Range('E6').Value = '  -4.47%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2959'
$ws.Range('E8').Value = '  -3.69%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07362'
$ws.Range('E9').Value = '  -4.31%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '22.79'
$ws.Range('E10').Value = '  -3.50%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07682'
$ws.Range('D12').Value = '1.832.59'
$ws.Range('E12').Value = '  -1.22%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.021'
$ws.Range('E13').Value = '  -2.75%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6753'
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '86.47'
$ws.Range('E15').Value = '  -5.28%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.187'
$ws.Range('E16').Value = '  -1.80%  '
$ws.Range('D17').Value = '29.170.96'
$ws.Range('E17').Value = '  -1.09%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008243'
$ws.Range('E18').Value = '  -1.16%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '229.17'
$ws.Range('E19').Value = '  -3.83%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.52'
$ws.Range('E20').Value = '  -1.78%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.308'
$ws.Range('E22').Value = '  -4.25%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.9992'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '160.92'
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1421'
$ws.Range('E25').Value = '  -4.99%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.678'
$ws.Range('E26').Value = '  -2.39%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.02'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.503'
$ws.Range('E28').Value = '  -1.79%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.233'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.099'
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.05334'
$ws.Range('E32').Value = '  +4.60%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.862'
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7456'
$ws.Range('E34').Value = '  -3.67%  '
$ws.Range('E35').Value = '  -1.80%  '
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('D37').Value = '1.315.45'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  -3.76%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.712'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9242'
$ws.Range('E40').Value = '  -2.92%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.005'
$ws.Range('E41').Value = '  +3.79%  '
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '103.55'
$ws.Range('E43').Value = '  -2.46%  '
$ws.Range('D44').Value = '1.986.52'
$ws.Range('E44').Value = '  -0.74%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.5165'
$ws.Range('E45').Value = '  -1.01%  '
$ws.Range('E46').Value = '  -3.16%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.762'
$ws.Range('E47').Value = '  -1.40%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '63.61'
$ws.Range('E48').Value = '  +0.69%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.314'
$ws.Range('E49').Value = '  -5.38%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07480'
$ws.Range('E50').Value = '  +9.80%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05926'
$ws.Range('E51').Value = '  -0.03%  '
